# Update "Use Case Beeldscherm OLC.docx"
$d = $word.ActiveDocument

# 1. Replace "TV" with "Software (voor beeldscherm)" in the Actoren row.
$d.Content.Find.Execute("Leerling, Schoolpas, TV", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Leerling, Schoolpas, Software (voor beeldscherm)", 2)

# 2. Remove the spell-check wavy-underline markers (proofErr) around "etc" in the Trigger row.
$d.Content.Find.Execute("(over school, huiswerk, etc)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "(over school, huiswerk, etc)", 2)
